# Fix 1.0.2 (Save as Type Into, Environment.CurrDir)
#
# Populates the process timestamps on row 2 ("SGP" operating unit) of
# Sheet1 with the date/time values recorded when that run completed,
# matching the style already used by the template rows below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the existing date/time number format (style index 1) from the
# template rows onto the cells we are about to fill in, so the new
# values render the same way as the rest of the sheet.
$ws.Range("B4:E4").Copy()
$ws.Range("B2:E2").PasteSpecial(-4122)

$ws.Range("G4:O4").Copy()
$ws.Range("G2:O2").PasteSpecial(-4122)

# Importation, Validation, Create Accounting, Transfer to GL
$ws.Range("B2").Value = 45625.785821759258
$ws.Range("C2").Value = 45625.788321759261
$ws.Range("D2").Value = 45625.790902777779
$ws.Range("E2").Value = 45625.795115740744

# Unaccounted, AP Trial Balance, Open Acct Balance, Invoice Aging,
# Payables Posted Invoice Register, Payables Posted Payment Register,
# Unpaid Invoices Report, Trial Balance Detail, Email
$ws.Range("G2").Value = 45625.816770833335
$ws.Range("H2").Value = 45625.818437499998
$ws.Range("I2").Value = 45625.819988425923
$ws.Range("J2").Value = 45625.821527777778
$ws.Range("K2").Value = 45625.823217592595
$ws.Range("L2").Value = 45625.824803240743
$ws.Range("M2").Value = 45625.826226851852
$ws.Range("N2").Value = 45625.828668981485
$ws.Range("O2").Value = 45625.829108796293
